# "Add graph search plots"
#
# The "Astar" sheet's results table (rows 2-62) previously had an extra block
# of UCS rows (CS1/CS2/Math x Low/Medium/High, rows 47-55) that used a
# placeholder load of 4,000,000 with no Avg/Time data. That block is removed,
# and the "Upper Bound" summary block (originally rows 56-62, one row per
# data set) moves up to follow directly after the existing data - but its
# "Physics" row (no matching data set in this sheet) is also dropped, since
# the table only tracks data sets that actually appear here.
#
# End result: rows 47-55 and the old "Physics" Upper Bound row are deleted,
# shrinking the table from A1:F62 down to A1:F52.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the placeholder UCS rows (CS1/CS2/Math x Low/Medium/High).
$ws.Range("A47:F55").EntireRow.Delete() | Out-Null

# Remove the "Upper Bound" / "Physics" row, which after the deletion above
# has shifted up from row 59 to row 50.
$ws.Range("A50:F50").EntireRow.Delete() | Out-Null

# Match the author's final selection: the whole of (new) row 46 selected,
# scrolled near the bottom of the now-shorter table.
$ws.Activate()
$ws.Rows.Item(46).Select() | Out-Null
